# requirements added, ontology starting to be clean
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Remove the three retired requirement rows ---------------------------
# (delete bottom-up so row indices of the still-to-delete rows don't shift)
# row 12 = speech-r11 "The comunication channels may have interactions"
# row 7  = speech-r6  "A political party can be identified"
# row 6  = speech-r5  "An election period can be identified"
$ws.Rows(12).Delete()
$ws.Rows(7).Delete()
$ws.Rows(6).Delete()

# After the deletions the sheet collapses to 25 data rows (rows 2-25) in the
# same relative order, just renumbered. We now touch up the rows whose text,
# style or height actually changed, and append the brand new row 26.

# --- 2. speech-r10 (row 9): reworded, font colour reset to automatic -------
# (these two rows used to carry the "red note" font; the reword drops that
# and the font goes back to the sheet's normal/automatic colour)
$ws.Range("B9").Value = "A political party may have a hashtag associated to a speech"
$ws.Range("B9").Font.ColorIndex = -4105
$ws.Range("B9").Font.Color = $ws.Range("A1").Font.Color
$ws.Range("B9").WrapText = $true
$ws.Rows(9).RowHeight = 34

# --- 3. speech-r12 (row 10): reworded, provenance note removed -------------
$ws.Range("B10").Value = "An speech may have interaction metrics"
$ws.Range("B10").Font.ColorIndex = -4105
$ws.Range("B10").Font.Color = $ws.Range("A1").Font.Color
$ws.Range("B10").WrapText = $true
$ws.Range("C10").ClearContents()
$ws.Rows(10).RowHeight = 17

# --- 4. speech-r15 (row 13): now holds real content + provenance note ------
# (provenance-backed rows use the red "note" font, same as the other
# speech-rX rows that carry a Column C comment)
$ws.Range("B13").Value = "Which are the threads previous to the 2021 elections in Madrid of all parties?"
$ws.Range("B13").Font.Color = 255
$ws.Range("B13").WrapText = $true
$ws.Range("C13").Value = "que hilo de qué?"
$ws.Rows(13).RowHeight = 34

# --- 5. speech-r16..speech-r27 (rows 14-25): fill in the new requirement text
$newRequirements = @{
    14 = "Which is the previous tweet of the elections of Madrid in 2015?"
    15 = "Which is the user most mentioned by a political party?"
    16 = "Which is the user most cited by a political party?"
    17 = "Which is the tweet with more citations of a political party?"
    18 = "Which are the tweets and users cited by at least two political parties?"
    19 = "Which are the tweets and users retweeted by at least two political parties?"
    20 = "Which is the most cited or retweeted tweet by the more number of political parties?"
    21 = "Which is the tweet of a political party most retweeted by another political party?"
    22 = "Which is the hashtag most tweeted by a political party in the 2019 Madrid elections?"
    23 = "Which is the most wetweeted hashtag by a political party in the 2019 Madrid elections?"
    24 = "Which is the political party that has tweeted the most in 2021?"
    25 = "Which is the political party that has been retweeted the most in 2021?"
}
foreach ($r in $newRequirements.Keys) {
    $ws.Range("B$r").Value = $newRequirements[$r]
    $ws.Range("B$r").WrapText = $true
    $ws.Rows($r).RowHeight = 34
}

# --- 6. brand new row 26: no identifier, just a requirement -----------------
$ws.Range("B26").Value = "How have ODS changed in the proposals of a manifesto of a political party?"
$ws.Range("B26").WrapText = $true
$ws.Rows(26).RowHeight = 34

# --- 7. selection / view bookkeeping ----------------------------------------
$ws.Range("B9").Select()
